$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row appended after the existing last row (row 13) -> row 14,
# matching the resale-numbers snapshot taken 2024-01-03 22:53:10.
$row = 14

# Columns A-D hold text (date / time / weekday / zero-padded week strings).
# Force text formatting before assignment so values like "2024-01-03" and
# "00" aren't auto-coerced into a date serial / number by the COM layer.
$textRange = $ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($row, 4))
$textRange.NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2024-01-03"
$ws.Cells.Item($row, 2).Value = "22:53:10"
$ws.Cells.Item($row, 3).Value = "Wednesday"
$ws.Cells.Item($row, 4).Value = "00"

# Strip the temporary formatting again so the new cells end up unstyled,
# same as the rest of this row's siblings (values stay text either way).
$textRange.ClearFormats()

# Columns E-T hold the numeric resale counts per city.
$ws.Cells.Item($row, 5).Value  = 140165
$ws.Cells.Item($row, 6).Value  = 142819
$ws.Cells.Item($row, 7).Value  = 172256
$ws.Cells.Item($row, 8).Value  = 146992
$ws.Cells.Item($row, 9).Value  = -1
$ws.Cells.Item($row, 10).Value = 117554
$ws.Cells.Item($row, 11).Value = 224003
$ws.Cells.Item($row, 12).Value = 248178
$ws.Cells.Item($row, 13).Value = 184195
$ws.Cells.Item($row, 14).Value = 109927
$ws.Cells.Item($row, 15).Value = 40233
$ws.Cells.Item($row, 16).Value = 30834
$ws.Cells.Item($row, 17).Value = 72233
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 41738
$ws.Cells.Item($row, 20).Value = -1
